# Generate Report for Handoff
# Update Priority ("low" -> "ht") and Latest Handoff Datetime for rows that
# just became ready for a fresh handoff, on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

$zhHandoffTime = "2016-08-30 10:32:32"
$deHandoffTime = "2016-08-30 10:32:37"

for ($r = 4; $r -le 7; $r++) {
    $zh.Cells.Item($r, 5).Value = "ht"
    $zh.Cells.Item($r, 8).Value = $zhHandoffTime

    $de.Cells.Item($r, 5).Value = "ht"
    $de.Cells.Item($r, 8).Value = $deHandoffTime

    # Overview sheet's "Latest HO Xliff Generate Date" column mirrors the
    # de-de handoff timestamp for these rows.
    $overview.Cells.Item($r, 7).Value = $deHandoffTime
}
